$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

$ws.Range("A2").Value = "AD.SEC.001.FON.02"
$ws.Range("A3").Value = "AD.SEC.001.FON.01"
$ws.Range("A4").Value = "AD.SEC.001.FON.03"
$ws.Range("A5").Value = "AD.DEP.001.FON.01"
$ws.Range("A6").Value = "RO.ACT"
$ws.Range("A7").Value = "RO.FOU"
$ws.Range("A8").Value = "MP.CPT"
$ws.Range("A9").Value = "RT.ART"
$ws.Range("A10").Value = "RT.MAT"
$ws.Range("A11").Value = "RO.ORG"
$ws.Range("A12").Value = "RT.EQU"
$ws.Range("A13").Value = "AD.SEC.014.FON.01"

$ws.Range("A4:A9").NumberFormat = "@"

$ws.Range("B12:B13").Select()
